$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E14").Value = "SWPA5040S101MT"
$ws.Range("F14").Value = "INDUCTOR SMD 100uF 0.75A 5x5mm"
$ws.Range("I14").Value = "C88132"

$ws.Range("K14").Select() | Out-Null
